$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value2 = 7009
$ws.Range("K3").Value2 = 7251
$ws.Range("I4").Value2 = 1809
$ws.Range("J4").Value2 = 1840
$ws.Range("K4").Value2 = 1499
$ws.Range("K6").Value2 = 7979
$ws.Range("I7").Value2 = 26268
$ws.Range("J7").Value2 = 29306
$ws.Range("K7").Value2 = 24254

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K5").Value2 = 64
$ws.Range("K8").Value2 = 1583
$ws.Range("K9").Value2 = 111
$ws.Range("K10").Value2 = 135
$ws.Range("K11").Value2 = 449
$ws.Range("K14").Value2 = 118
$ws.Range("K15").Value2 = 252
$ws.Range("K18").Value2 = 163
$ws.Range("K19").Value2 = 714
$ws.Range("K20").Value2 = 590
$ws.Range("K23").Value2 = 243
$ws.Range("K24").Value2 = 76
$ws.Range("K25").Value2 = 111
$ws.Range("K29").Value2 = 1332
$ws.Range("K31").Value2 = 272
$ws.Range("K33").Value2 = 1040
$ws.Range("K37").Value2 = 819
$ws.Range("K41").Value2 = 167
$ws.Range("K42").Value2 = 893
$ws.Range("K43").Value2 = 196
$ws.Range("K44").Value2 = 201
$ws.Range("K47").Value2 = 164
$ws.Range("K48").Value2 = 313
$ws.Range("K51").Value2 = 307
$ws.Range("K52").Value2 = 635
$ws.Range("K53").Value2 = 307
$ws.Range("K54").Value2 = 472
$ws.Range("K55").Value2 = 263
$ws.Range("K59").Value2 = 43
$ws.Range("I63").Value2 = 230
$ws.Range("J63").Value2 = 120
$ws.Range("K63").Value2 = 69
$ws.Range("K67").Value2 = 946
$ws.Range("K71").Value2 = 75
$ws.Range("K73").Value2 = 218
$ws.Range("K76").Value2 = 325
$ws.Range("K78").Value2 = 291
$ws.Range("K79").Value2 = 592
$ws.Range("K80").Value2 = 87
$ws.Range("K83").Value2 = 518
$ws.Range("K85").Value2 = 1115
$ws.Range("K88").Value2 = 261
$ws.Range("K89").Value2 = 360
$ws.Range("K94").Value2 = 324
$ws.Range("K96").Value2 = 262
$ws.Range("K98").Value2 = 125
$ws.Range("I101").Value2 = 26268
$ws.Range("J101").Value2 = 29306
$ws.Range("K101").Value2 = 24254

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K3").Value2 = 30
$ws.Range("K7").Value2 = 118

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value2 = 111
$ws.Range("K7").Value2 = 262

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K6").Value2 = 155
$ws.Range("K7").Value2 = 449

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K6").Value2 = 105
$ws.Range("K7").Value2 = 360

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value2 = 367
$ws.Range("K6").Value2 = 273
$ws.Range("K7").Value2 = 1115

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value2 = 232
$ws.Range("K7").Value2 = 635

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value2 = 79
$ws.Range("K7").Value2 = 307

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K3").Value2 = 481
$ws.Range("K7").Value2 = 1583

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value2 = 180
$ws.Range("K7").Value2 = 518

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value2 = 260
$ws.Range("K3").Value2 = 371
$ws.Range("K4").Value2 = 53
$ws.Range("K7").Value2 = 1040

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value2 = 237
$ws.Range("K3").Value2 = 269
$ws.Range("K6").Value2 = 243
$ws.Range("K7").Value2 = 819

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value2 = 101
$ws.Range("K7").Value2 = 272

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value2 = 259
$ws.Range("K3").Value2 = 342
$ws.Range("K7").Value2 = 946

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K4").Value2 = 28
$ws.Range("K7").Value2 = 472

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value2 = 377
$ws.Range("K3").Value2 = 474
$ws.Range("K6").Value2 = 389
$ws.Range("K7").Value2 = 1332

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value2 = 75
$ws.Range("K7").Value2 = 313

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value2 = 208
$ws.Range("K3").Value2 = 214
$ws.Range("K6").Value2 = 237
$ws.Range("K7").Value2 = 714

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K2").Value2 = 55
$ws.Range("K7").Value2 = 201

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K3").Value2 = 62
$ws.Range("K7").Value2 = 325

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K4").Value2 = 10
$ws.Range("K7").Value2 = 167

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K4").Value2 = 38
$ws.Range("K6").Value2 = 334
$ws.Range("K7").Value2 = 893

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K3").Value2 = 23
$ws.Range("K7").Value2 = 135

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value2 = 86
$ws.Range("K6").Value2 = 98
$ws.Range("K7").Value2 = 291

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K3").Value2 = 77
$ws.Range("K7").Value2 = 263

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K3").Value2 = 19
$ws.Range("K7").Value2 = 76

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value2 = 84
$ws.Range("K6").Value2 = 66
$ws.Range("K7").Value2 = 243

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value2 = 199
$ws.Range("K7").Value2 = 592

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value2 = 203
$ws.Range("K3").Value2 = 190
$ws.Range("K6").Value2 = 163
$ws.Range("K7").Value2 = 590

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K6").Value2 = 43
$ws.Range("K7").Value2 = 163

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K3").Value2 = 68
$ws.Range("K6").Value2 = 148
$ws.Range("K7").Value2 = 324

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K2").Value2 = 43
$ws.Range("K7").Value2 = 111

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K2").Value2 = 47
$ws.Range("K3").Value2 = 49
$ws.Range("K7").Value2 = 164

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value2 = 95
$ws.Range("K7").Value2 = 252

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K2").Value2 = 22
$ws.Range("K6").Value2 = 73
$ws.Range("K7").Value2 = 125

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K4").Value2 = 7
$ws.Range("K7").Value2 = 111

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K6").Value2 = 74
$ws.Range("K7").Value2 = 218

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("K3").Value2 = 13
$ws.Range("K7").Value2 = 43

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K6").Value2 = 104
$ws.Range("K7").Value2 = 261

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K6").Value2 = 30
$ws.Range("K7").Value2 = 64

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K3").Value2 = 84
$ws.Range("K4").Value2 = 33
$ws.Range("K6").Value2 = 101
$ws.Range("K7").Value2 = 307

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K6").Value2 = 74
$ws.Range("K7").Value2 = 196

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K6").Value2 = 20
$ws.Range("K7").Value2 = 75

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K3").Value2 = 18
$ws.Range("K7").Value2 = 87
